$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.850.19'
$ws.Range("E2").Value = '  -2.32%  '

$ws.Range("D3").Value = '1.815.03'
$ws.Range("E3").Value = '  -1.36%  '

$ws.Range("E4").Value = '  -0.53%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.008'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.30'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4610'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3639'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07217'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8550'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.71'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.91%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.845.56'
$ws.Range("E12").Value = '  -2.24%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07537'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.321'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.53%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.75'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.62%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.485'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008592'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.008'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.41%  '

$ws.Range("D20").Value = '27.005.13'
$ws.Range("E20").Value = '  -1.80%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.41'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.143'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.49'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.08%  '

$ws.Range("D24").Value = '2.068.81'
$ws.Range("E24").Value = '  -1.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.56'
$ws.Range("D25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.847'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.57%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.12'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.063'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.077'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.55'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08855'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.961'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.129'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.51%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.392'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.73%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7184'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.008'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.55%  '

$ws.Range("E37").Value = '  -3.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05226'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.416'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01913'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.905'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.126'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5131'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1617'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.166'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4781'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.07%  '

$ws.Range("E47").Value = '  -0.51%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.00'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.84%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.07'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.612'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06197'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.98%  '
